# Auto-generated Excel COM-interop script to apply the Diabolos_Profits update.
# Updates currentAveragePrice* / Leve profit columns (H-N) for specific leves
# across all 8 sheets, matching the "update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 706.5294
$ws.Range("I17").Value = 880.5
$ws.Range("J17").Value = 653
$ws.Range("K17").Value = 2641.5
$ws.Range("L17").Value = 1959
$ws.Range("M17").Value = -2473.5
$ws.Range("N17").Value = -2295
# Row 33
$ws.Range("H33").Value = 999.7222
$ws.Range("I33").Value = 1067.2142
$ws.Range("K33").Value = 1067.2142
$ws.Range("M33").Value = -838.2141999999999
# Row 53
$ws.Range("H53").Value = 2765
$ws.Range("I53").Value = 70.583336
$ws.Range("K53").Value = 70.583336
$ws.Range("M53").Value = 566.416664
# Row 62
$ws.Range("H62").Value = 987218.5600000001
$ws.Range("I62").Value = 2063880.8
$ws.Range("J62").Value = 90000
$ws.Range("K62").Value = 2063880.8
$ws.Range("L62").Value = 90000
$ws.Range("M62").Value = -2063256.8
$ws.Range("N62").Value = -91248
# Row 65
$ws.Range("H65").Value = 987218.5600000001
$ws.Range("I65").Value = 2063880.8
$ws.Range("J65").Value = 90000
$ws.Range("K65").Value = 10319404
$ws.Range("L65").Value = 450000
$ws.Range("M65").Value = -10316284
$ws.Range("N65").Value = -456240
# Row 69
$ws.Range("H69").Value = 10420.2
$ws.Range("J69").Value = 10420.2
$ws.Range("L69").Value = 31260.6
$ws.Range("N69").Value = -33008.60000000001
# Row 72
$ws.Range("H72").Value = 10420.2
$ws.Range("J72").Value = 10420.2
$ws.Range("L72").Value = 93781.8
$ws.Range("N72").Value = -102517.8
# Row 111
$ws.Range("H111").Value = 28511.916
$ws.Range("J111").Value = 78354.336
$ws.Range("L111").Value = 235063.008
$ws.Range("N111").Value = -241197.008
# Row 137
$ws.Range("H137").Value = 15627369
$ws.Range("I137").Value = 50001960
$ws.Range("J137").Value = 2555.682
$ws.Range("K137").Value = 150005880
$ws.Range("L137").Value = 7667.045999999999
$ws.Range("M137").Value = -150003330
$ws.Range("N137").Value = -12767.046
# Row 138
$ws.Range("H138").Value = 3178.3928
$ws.Range("I138").Value = 2194.2307
$ws.Range("J138").Value = 4031.3333
$ws.Range("K138").Value = 6582.6921
$ws.Range("L138").Value = 12093.9999
$ws.Range("M138").Value = -1442.6921
$ws.Range("N138").Value = -22373.9999
# Row 141
$ws.Range("H141").Value = 2380.8096
$ws.Range("I141").Value = 2247
$ws.Range("J141").Value = 2949.5
$ws.Range("K141").Value = 6741
$ws.Range("L141").Value = 8848.5
$ws.Range("M141").Value = -1561
$ws.Range("N141").Value = -19208.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 424.11
$ws.Range("I32").Value = 385.39584
$ws.Range("J32").Value = 1353.25
$ws.Range("K32").Value = 385.39584
$ws.Range("L32").Value = 1353.25
$ws.Range("M32").Value = -98.39584000000002
$ws.Range("N32").Value = -1927.25
# Row 61
$ws.Range("H61").Value = 2146.2144
$ws.Range("I61").Value = 1060.8889
$ws.Range("J61").Value = 4099.8
$ws.Range("K61").Value = 1060.8889
$ws.Range("L61").Value = 4099.8
$ws.Range("M61").Value = -848.8888999999999
$ws.Range("N61").Value = -4523.8
# Row 122
$ws.Range("H122").Value = 4143.6924
$ws.Range("I122").Value = 3525.9333
$ws.Range("K122").Value = 10577.7999
$ws.Range("M122").Value = -8127.7999
# Row 136
$ws.Range("H136").Value = 2146.2144
$ws.Range("I136").Value = 1060.8889
$ws.Range("J136").Value = 4099.8
$ws.Range("K136").Value = 3182.6667
$ws.Range("L136").Value = 12299.4
$ws.Range("M136").Value = -632.6666999999998
$ws.Range("N136").Value = -17399.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 15
$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 500
$ws.Range("K15").Value = 500
$ws.Range("M15").Value = -273
# Row 82
$ws.Range("H82").Value = 37713.855
$ws.Range("I82").Value = 9500.5
$ws.Range("J82").Value = 48999.2
$ws.Range("K82").Value = 9500.5
$ws.Range("L82").Value = 48999.2
$ws.Range("M82").Value = -9117.5
$ws.Range("N82").Value = -49765.2
# Row 85
$ws.Range("H85").Value = 37713.855
$ws.Range("I85").Value = 9500.5
$ws.Range("J85").Value = 48999.2
$ws.Range("K85").Value = 9500.5
$ws.Range("L85").Value = 48999.2
$ws.Range("M85").Value = -8174.5
$ws.Range("N85").Value = -51651.2
# Row 134
$ws.Range("H134").Value = 2011.4762
$ws.Range("I134").Value = 1737.05
$ws.Range("K134").Value = 5211.15
$ws.Range("M134").Value = -2676.15

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1200
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -4224
# Row 62
$ws.Range("H62").Value = 6313.0625
$ws.Range("I62").Value = 6193.077
$ws.Range("J62").Value = 6833
$ws.Range("K62").Value = 6193.077
$ws.Range("L62").Value = 6833
$ws.Range("M62").Value = -5569.077
$ws.Range("N62").Value = -8081
# Row 65
$ws.Range("H65").Value = 6313.0625
$ws.Range("I65").Value = 6193.077
$ws.Range("J65").Value = 6833
$ws.Range("K65").Value = 30965.385
$ws.Range("L65").Value = 34165
$ws.Range("M65").Value = -27845.385
$ws.Range("N65").Value = -40405

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 9412.727999999999
$ws.Range("I6").Value = 12566.75
$ws.Range("J6").Value = 1002
$ws.Range("K6").Value = 37700.25
$ws.Range("L6").Value = 3006
$ws.Range("M6").Value = -37587.25
$ws.Range("N6").Value = -3232
# Row 131
$ws.Range("H131").Value = 12609.789
$ws.Range("I131").Value = 718.6667
$ws.Range("J131").Value = 23311.8
$ws.Range("K131").Value = 2156.0001
$ws.Range("L131").Value = 69935.39999999999
$ws.Range("M131").Value = 2883.9999
$ws.Range("N131").Value = -80015.39999999999
# Row 134
$ws.Range("H134").Value = 1840.2858
$ws.Range("I134").Value = 1813.6666
$ws.Range("K134").Value = 5440.9998
$ws.Range("M134").Value = -370.9997999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2395.9092
$ws.Range("I80").Value = 1371.8
$ws.Range("K80").Value = 1371.8
$ws.Range("M80").Value = -373.8
# Row 83
$ws.Range("H83").Value = 2395.9092
$ws.Range("I83").Value = 1371.8
$ws.Range("K83").Value = 6859
$ws.Range("M83").Value = -1867
# Row 95
$ws.Range("H95").Value = 58000
$ws.Range("J95").Value = 58000
$ws.Range("L95").Value = 58000
$ws.Range("N95").Value = -63492
# Row 97
$ws.Range("H97").Value = 845.3158
$ws.Range("J97").Value = 629.3
$ws.Range("L97").Value = 629.3
$ws.Range("N97").Value = -1621.3
# Row 102
$ws.Range("H102").Value = 1943.3243
$ws.Range("I102").Value = 1382.3928
$ws.Range("K102").Value = 1382.3928
$ws.Range("M102").Value = 239.6071999999999
# Row 122
$ws.Range("H122").Value = 2865.1177
$ws.Range("I122").Value = 2408.3333
$ws.Range("J122").Value = 3961.4
$ws.Range("K122").Value = 7224.999899999999
$ws.Range("L122").Value = 11884.2
$ws.Range("M122").Value = -4774.999899999999
$ws.Range("N122").Value = -16784.2
# Row 136
$ws.Range("H136").Value = 65000
$ws.Range("J136").Value = 65000
$ws.Range("L136").Value = 195000
$ws.Range("N136").Value = -200100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 150
$ws.Range("J14").Value = 150
$ws.Range("L14").Value = 150
$ws.Range("N14").Value = -494
# Row 16
$ws.Range("H16").Value = 2731.2917
$ws.Range("I16").Value = 2235.3845
$ws.Range("J16").Value = 3317.3635
$ws.Range("K16").Value = 2235.3845
$ws.Range("L16").Value = 3317.3635
$ws.Range("M16").Value = -2065.3845
$ws.Range("N16").Value = -3657.3635
# Row 61
$ws.Range("H61").Value = 6391.6562
$ws.Range("I61").Value = 7258.96
$ws.Range("J61").Value = 3294.1428
$ws.Range("K61").Value = 7258.96
$ws.Range("L61").Value = 3294.1428
$ws.Range("M61").Value = -7056.96
$ws.Range("N61").Value = -3698.1428
# Row 113
$ws.Range("H113").Value = 6391.6562
$ws.Range("I113").Value = 7258.96
$ws.Range("J113").Value = 3294.1428
$ws.Range("K113").Value = 7258.96
$ws.Range("L113").Value = 3294.1428
$ws.Range("M113").Value = -5088.96
$ws.Range("N113").Value = -7634.1428
# Row 122
$ws.Range("H122").Value = 3870.5557
$ws.Range("I122").Value = 3015.75
$ws.Range("K122").Value = 9047.25
$ws.Range("M122").Value = -6597.25
# Row 133
$ws.Range("H133").Value = 39995
$ws.Range("J133").Value = 39995
$ws.Range("L133").Value = 39995
$ws.Range("N133").Value = -45055

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 959076.9399999999
$ws.Range("I62").Value = 2650614
$ws.Range("J62").Value = 7587.3125
$ws.Range("K62").Value = 2650614
$ws.Range("L62").Value = 7587.3125
$ws.Range("M62").Value = -2649990
$ws.Range("N62").Value = -8835.3125
# Row 65
$ws.Range("H65").Value = 959076.9399999999
$ws.Range("I65").Value = 2650614
$ws.Range("J65").Value = 7587.3125
$ws.Range("K65").Value = 13253070
$ws.Range("L65").Value = 37936.5625
$ws.Range("M65").Value = -13249950
$ws.Range("N65").Value = -44176.5625
# Row 70
$ws.Range("H70").Value = 40000
$ws.Range("I70").Value = 35000
$ws.Range("K70").Value = 35000
$ws.Range("M70").Value = -34685
# Row 73
$ws.Range("H73").Value = 40000
$ws.Range("I73").Value = 35000
$ws.Range("K73").Value = 35000
$ws.Range("M73").Value = -33908
# Row 81
$ws.Range("H81").Value = 22228976
$ws.Range("I81").Value = 1199
$ws.Range("J81").Value = 25007448
$ws.Range("K81").Value = 2398
$ws.Range("L81").Value = 50014896
$ws.Range("M81").Value = -1337
$ws.Range("N81").Value = -50017018
# Row 84
$ws.Range("H84").Value = 22228976
$ws.Range("I84").Value = 1199
$ws.Range("J84").Value = 25007448
$ws.Range("K84").Value = 11990
$ws.Range("L84").Value = 250074480
$ws.Range("M84").Value = -6686
$ws.Range("N84").Value = -250085088
# Row 96
$ws.Range("H96").Value = 3649.7334
$ws.Range("I96").Value = 2495.182
$ws.Range("J96").Value = 6824.75
$ws.Range("K96").Value = 2495.182
$ws.Range("L96").Value = 6824.75
$ws.Range("M96").Value = -1122.182
$ws.Range("N96").Value = -9570.75
# Row 107
$ws.Range("H107").Value = 426.46667
$ws.Range("I107").Value = 407.83334
$ws.Range("J107").Value = 501
$ws.Range("K107").Value = 1223.50002
$ws.Range("L107").Value = 1503
$ws.Range("M107").Value = 696.4999800000001
$ws.Range("N107").Value = -5343
# Row 122
$ws.Range("H122").Value = 2699.95
$ws.Range("I122").Value = 1623.4706
$ws.Range("J122").Value = 8800
$ws.Range("K122").Value = 4870.4118
$ws.Range("L122").Value = 26400
$ws.Range("M122").Value = -2420.4118
$ws.Range("N122").Value = -31300
# Row 136
$ws.Range("H136").Value = 4218.8335
$ws.Range("I136").Value = 3283.0952
$ws.Range("K136").Value = 9849.285600000001
$ws.Range("M136").Value = -7299.285600000001
